$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "DFejViks2uY,DFejViks2uY,DFejViks2uY,RR9qtpQL0mY,cElhIDdGz7M,PL2eqHfPK3Q,cElhIDdGz7M,V7B9B_D0jjg,cElhIDdGz7M,SPSJGSLdFHc,cElhIDdGz7M,BU619veKVHQ,cElhIDdGz7M,F_R3uGDXhAc,cElhIDdGz7M,Lf1QfapEL7c,cElhIDdGz7M,LdZZKGXf5G4,cElhIDdGz7M,F_R3uGDXhAc,cElhIDdGz7M,joK0ZhaHeO4,cElhIDdGz7M,V7B9B_D0jjg,cElhIDdGz7M,SPSJGSLdFHc,cElhIDdGz7M,BU619veKVHQ,cElhIDdGz7M,8u4tcWt8uew,cElhIDdGz7M,F_R3uGDXhAc,cElhIDdGz7M,iyyVoLv513E,cElhIDdGz7M,dM4Ugwy2tRQ,cElhIDdGz7M,GUQKIROlp2w,cElhIDdGz7M,1qt6ijv1T6I,cElhIDdGz7M,9EC8yodm65Q,cElhIDdGz7M,tSZ6qC7GcdA,cElhIDdGz7M,F_R3uGDXhAc,cElhIDdGz7M,Lf1QfapEL7c,cElhIDdGz7M,KTujQiqqtcs,cElhIDdGz7M,CDQsRe6GZHU,cElhIDdGz7M,0iOheWSY2ng,cElhIDdGz7M,zbEJ8ljMev0,cElhIDdGz7M,c4tlzZrBses,cElhIDdGz7M,moZHzH71Zxk,cElhIDdGz7M,PT39Dl687IM,cElhIDdGz7M,MQ4wHdDUqkQ,cElhIDdGz7M,LSZIjPaZTb4,cElhIDdGz7M,DFejViks2uY,cElhIDdGz7M,9Zs2CJTYr48,cElhIDdGz7M,50m5WnlzsI8,cElhIDdGz7M,5_N7bzip72k,cElhIDdGz7M,DFejViks2uY,cElhIDdGz7M,0_sv895JPlc,0_sv895JPlc,0_sv895JPlc,sQkDUCblDUo,cElhIDdGz7M,LSZIjPaZTb4,cElhIDdGz7M,u41ujNodvnM,cElhIDdGz7M,h5rmV9bvryQ,cElhIDdGz7M,SD_Qd4q3Vlo,cElhIDdGz7M,8u4tcWt8uew,cElhIDdGz7M,zLWCootfQ64,cElhIDdGz7M,dM4Ugwy2tRQ,cElhIDdGz7M,GogAr_2DMWc,cElhIDdGz7M,D1cmuQT3ENA,cElhIDdGz7M,OVx4IgI3rgg,cElhIDdGz7M,dM4Ugwy2tRQ,cElhIDdGz7M,6TgbuYkiEIs,cElhIDdGz7M,8u4tcWt8uew,cElhIDdGz7M,dM4Ugwy2tRQ,cElhIDdGz7M,6TgbuYkiEIs,cElhIDdGz7M,5tOuYW7k9LM,cElhIDdGz7M,GogAr_2DMWc,cElhIDdGz7M,s61VYoAdp24,cElhIDdGz7M,8ElG31cR0Ng,cElhIDdGz7M"
$ws.Range("B2").Value = "261,261,261,1112,5,2224,5,452,5,177,5,3073,5,22,5,61,5,209,5,22,5,1898,5,452,5,177,5,3073,5,817,5,22,5,2553,5,953,5,1597,5,135,5,474,5,504,5,22,5,61,5,733,5,314,5,939,5,98,5,590,5,1027,5,5250,5,1775,5,5907,5,260,5,668,5,4340,5,120,5,261,5,157,157,157,33,5,3004,5,1666,5,1621,5,939,5,817,5,161,5,953,5,607,5,142,5,343,5,953,5,1059,5,817,5,953,5,1059,5,2746,5,607,5,3569,5,2481,5"
$ws.Range("C2").Value = "266,266,266,1113,7,2225,7,454,8,179,8,3076,9,23,7,63,8,210,7,23,7,1900,8,454,8,179,8,3076,9,818,7,23,7,2554,7,955,8,1598,7,136,7,476,8,506,8,23,7,63,8,735,8,316,8,941,8,100,8,592,8,1029,8,5251,7,1776,7,5909,8,262,8,670,8,4342,8,122,8,266,11,163,163,163,35,8,3005,7,1668,8,1623,8,940,7,818,7,163,8,955,8,609,8,143,7,344,7,955,8,1061,8,818,7,955,8,1061,8,2747,7,609,8,3571,8,2483,8"
